$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Discontinued products to remove from the inventory table.
$namesToRemove = @(
    "ATROVENT 250MCG/2ML 20 UNIT DOSE VIAL",
    "FARCOLIN RESPIRATOR 0.5% SOLN. 20 ML",
    "PULMICORT 0.25MG/ML 20 NEBULIZER VIAL SUSP."
)

# The data rows of the table run from row 4 down to row 27 (24 products).
# Row 28 holds the grand-total row and row 29 the footer row (before any edits,
# these currently live at rows 31 and 32 because the table still has 27 products).
$firstDataRow = 4
$lastDataRow = 27 + $namesToRemove.Count   # 30 before trimming

foreach ($name in $namesToRemove) {
    $found = $ws.Cells.Find($name)
    if ($found -eq $null) { continue }
    $targetRow = $found.Row

    # Shift the content of the merged "item" columns (B:G, H:K, L:M and N) up by
    # one row, for every row below the removed item down to the last data row.
    # Column A (the running counter) is intentionally left untouched.
    for ($r = $targetRow; $r -lt $lastDataRow; $r++) {
        $srcRow = $r + 1
        $ws.Range("B$r").Value2 = $ws.Range("B$srcRow").Value2
        $ws.Range("H$r").Value2 = $ws.Range("H$srcRow").Value2
        $ws.Range("L$r").Value2 = $ws.Range("L$srcRow").Value2
        $ws.Range("N$r").Value2 = $ws.Range("N$srcRow").Value2
    }
}

# The trailing rows now hold stale, duplicated data (as many as products removed).
# Delete them outright; this also pulls the grand-total and footer rows back up.
for ($i = 0; $i -lt $namesToRemove.Count; $i++) {
    $ws.Rows.Item($lastDataRow - $i).Delete()
}
